$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Иванов "
$ws.Range("B1").Value = "Васильков"
$ws.Range("A2").Value = "Петренко"
$ws.Range("B2").Value = "Барахты"
$ws.Range("A3").Value = "Махно "
$ws.Range("B3").Value = "Козятие"
$ws.Range("A4").Value = "Степан"
$ws.Range("B4").Value = "Жорновка "
